$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix row ordering bug in onStepEnd: the reel-strip weighting table (rows 2-21)
# was being written in the wrong row order. Rewrite each row with the correct values.
$rowData = @{
    2 = @(1203, 3, 15, 15, 15, 15)
    3 = @(401, 9, 48, 67, 75, 45)
    4 = @(601, 9, 60, 67, 60, 42)
    5 = @(902, 1, 0, 0, 0, 0)
    6 = @(1201, 2, 10, 10, 10, 10)
    7 = @(101, 9, 30, 15, 60, 15)
    8 = @(901, 16, 15, 45, 60, 60)
    9 = @(801, 3, 67, 65, 52, 45)
    10 = @(1202, 2, 10, 10, 10, 10)
    11 = @(1001, 18, 30, 75, 60, 72)
    12 = @(701, 3, 90, 45, 97, 15)
    13 = @(201, 9, 30, 15, 45, 30)
    14 = @(301, 6, 45, 30, 60, 45)
    15 = @(501, 9, 52, 30, 75, 45)
    16 = @(1101, 0, 15, 30, 30, 0)
    17 = @(1, 0, 2, 2, 2, 2)
    18 = @(2, 0, 2, 2, 2, 2)
    19 = @(3, 0, 3, 3, 3, 3)
    20 = @(502, 0, 4, 0, 0, 0)
    21 = @(802, 0, 4, 5, 4, 0)
}

foreach ($r in $rowData.Keys) {
    $vals = $rowData[$r]
    for ($c = 1; $c -le 6; $c++) {
        $ws.Cells.Item($r, $c).Value = $vals[$c - 1]
    }
}

Write-Output "Updated rows 2-21 with corrected reel weight ordering"
